$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Assigned but not busy"
